$wb = $excel.ActiveWorkbook

# Update visitor/attendance counts on the "展览" (Exhibitions) sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 7632   # 南宁·2024三月三国潮动漫节（良牙春典）
$wsExhibit.Range("F5").Value = 34     # 南宁·第五人格Only1.0
$wsExhibit.Range("F6").Value = 276    # 南宁·AP动漫游戏嘉年华
$wsExhibit.Range("F10").Value = 152   # 南宁·AB动漫游戏嘉年华

# Update the same rows (shifted by one due to an extra "演出" row) on the
# "全部类型" (All types) combined sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7632    # 南宁·2024三月三国潮动漫节（良牙春典）
$wsAll.Range("F5").Value = 34      # 南宁·第五人格Only1.0
$wsAll.Range("F6").Value = 276     # 南宁·AP动漫游戏嘉年华
$wsAll.Range("F11").Value = 152    # 南宁·AB动漫游戏嘉年华
